# Adds a new data-row (row 17) to the end of the "Artfynd" sheet, matching
# the shape of the existing species-observation rows above it.
#
# Text cells are written with a leading apostrophe so that numeric- or
# date-looking strings ("2", "2023-08-29", ...) are kept as literal text
# instead of being auto-converted to a number/date by the host, then the
# cell style is reset to "Normal" so no incidental number-format sticks to
# the cell (keeping it equivalent to the untouched cells around it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 17

function Set-TextCell($rowIdx, $col, $text) {
    $ws.Cells.Item($rowIdx, $col).Value() = "'" + $text
    $ws.Cells.Item($rowIdx, $col).Style = "Normal"
}

function Set-NumberCell($rowIdx, $col, $num) {
    $ws.Cells.Item($rowIdx, $col).Value() = $num
}

function Set-BoolCell($rowIdx, $col, $flag) {
    $ws.Cells.Item($rowIdx, $col).Value() = $flag
}

# A Id
Set-NumberCell $row 1 111961716
# B Taxonsorteringsordning
Set-NumberCell $row 2 81076
# C Valideringsstatus
Set-TextCell $row 3 "Ovaliderad"
# D Rödlistade
Set-TextCell $row 4 "LC"
# E TaxonId
Set-NumberCell $row 5 5046
# F Artnamn
Set-TextCell $row 6 "Grön jordtunga"
# G Vetenskapligt namn
Set-TextCell $row 7 "Microglossum viride"
# H Auktor
Set-TextCell $row 8 "(Pers.:Fr.) Gillet"
# I Antal
Set-TextCell $row 9 "2"
# J Enhet
Set-TextCell $row 10 "mycel"
# K Ålder-Stadium (empty)
Set-TextCell $row 11 ""
# N Metod (empty)
Set-TextCell $row 14 ""
# P Lokalnamn
Set-TextCell $row 16 "Klockarbäcken, Vb"
# Q Ost
Set-NumberCell $row 17 753078.7913326195
# R Nord
Set-NumberCell $row 18 7090973.389402887
# S Noggrannhet
Set-NumberCell $row 19 100
# T Län
Set-TextCell $row 20 "Västerbotten"
# U Kommun
Set-TextCell $row 21 "Umeå"
# V Provins
Set-TextCell $row 22 "Västerbotten"
# W Församling
Set-TextCell $row 23 "Umeå socken"
# Y Startdatum
Set-TextCell $row 25 "2023-08-29"
# Z Starttid
Set-TextCell $row 26 "00:00"
# AA Slutdatum
Set-TextCell $row 27 "2023-08-29"
# AB Sluttid
Set-TextCell $row 28 "00:00"
# AD Ej återfunnen
Set-BoolCell $row 30 $false
# AE Osäker artbestämning
Set-BoolCell $row 31 $false
# AF Bestämningsmetod
Set-TextCell $row 32 "mikroskoperad"
# AG Ospontan
Set-BoolCell $row 33 $false
# AT Bestämningsår (empty)
Set-TextCell $row 46 ""
# AW Rapportör
Set-TextCell $row 49 "Stefan Phalagorn Bergström"
# AX Observatörer
Set-TextCell $row 50 "Stefan Phalagorn Bergström, Andreas Estensen, Annika  Carlberg , Ola Elleström, Thomas Strid, Anne Järvinen, Emma Sewell"
# AY Projektnamn (empty)
Set-TextCell $row 51 ""
